# Update "想去人数" (interested-count) figures on the 展览 (Worksheets 1)
# and 全部类型 (Worksheets 4) sheets to match the newly generated data
# (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 1056
$ws1.Range("F5").Value  = 13907
$ws1.Range("F6").Value  = 46
$ws1.Range("F8").Value  = 175
$ws1.Range("F11").Value = 141
$ws1.Range("F12").Value = 93
$ws1.Range("F14").Value = 536
$ws1.Range("F18").Value = 13961
$ws1.Range("F20").Value = 628
$ws1.Range("F21").Value = 14977
$ws1.Range("F23").Value = 8276
$ws1.Range("F24").Value = 276
$ws1.Range("F26").Value = 26
$ws1.Range("F28").Value = 425
$ws1.Range("F30").Value = 1
$ws1.Range("F43").Value = 392
$ws1.Range("F45").Value = 5095

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 1056
$ws4.Range("F5").Value  = 13907
$ws4.Range("F6").Value  = 46
$ws4.Range("F8").Value  = 176
$ws4.Range("F11").Value = 141
$ws4.Range("F12").Value = 93
$ws4.Range("F14").Value = 536
$ws4.Range("F18").Value = 13961
$ws4.Range("F20").Value = 628
$ws4.Range("F21").Value = 14977
$ws4.Range("F23").Value = 8276
$ws4.Range("F24").Value = 276
$ws4.Range("F26").Value = 26
$ws4.Range("F28").Value = 425
$ws4.Range("F30").Value = 1
$ws4.Range("F45").Value = 392
$ws4.Range("F47").Value = 5095
